$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.068.54'
$ws.Range('E2').Value = '  +6.55%  '
$ws.Range('D3').Value = '3.106.12'
$ws.Range('E3').Value = '  +4.66%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '584.28'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +4.14%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '144.53'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +5.48%  '
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('D8').Value = '3.099.19'
$ws.Range('E8').Value = '  +4.56%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.529'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +1.89%  '
$ws.Range('E10').Value = '  +13.73%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.76'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +8.09%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.466'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +3.19%  '
$ws.Range('E13').Value = '  +8.41%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '35.48'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +5.69%  '
$ws.Range('E15').Value = '  +0.70%  '
$ws.Range('D16').Value = '3.618.18'
$ws.Range('E16').Value = '  +4.65%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '7.15'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +1.28%  '
$ws.Range('D18').Value = '62.995.68'
$ws.Range('E18').Value = '  +6.38%  '
$ws.Range('D19').Value = '3.100.86'
$ws.Range('E19').Value = '  +4.43%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '464.23'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +6.85%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.16'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +4.92%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.726'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +1.30%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.51'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +7.38%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '13.28'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.07%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '81.80'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +2.71%  '
$ws.Range('E26').Value = '  -0.02%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.44'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +9.64%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.23'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.23%  '
$ws.Range('B29').Value = 'PancakeSwap'
$ws.Range('C29').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.67'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +5.42%  '
$ws.Range('B30').Value = 'FirstDigitalUSD'
$ws.Range('C30').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.00'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.05%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.83'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +10.27%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '26.90'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +5.04%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.111'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +4.56%  '
$ws.Range('D34').Value = '0.0₃0865'
$ws.Range('E34').Value = '  +13.78%  '
$ws.Range('E35').Value = '  +16.72%  '
$ws.Range('E36').Value = '  +5.52%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.33'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +21.48%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '6.04'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +2.95%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '50.80'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +4.72%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '433.46'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +8.83%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '8.72'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.81%  '
$ws.Range('D42').Value = '2.905.49'
$ws.Range('E42').Value = '  +5.98%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0369'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +5.05%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.281'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +12.75%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.112'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +7.15%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.16'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +8.17%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '35.10'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +2.71%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '122.73'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.32%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '24.48'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +5.31%  '
